# Apply the "Actualización automática" data refresh to the three sheets.
$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Sheet "VENTAS POR GRUPO": PORCELANATO sales for HIDALGO HIDALGO PEDRO GUSTAVO / CARRION CARRION LESLY ANABE
$wsVentasGrupo.Range("M5").Value = 4994.11

# Sheet "VENTA MENSUAL": julio sales for the same advisor/client, plus the julio column total
$wsVentaMensual.Range("F5").Value = 5395.92
$wsVentaMensual.Range("F22").Value = 37623.7

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO group row (VENTA / POR CUMPLIR / CUMPLIMIENTO)
$wsCumplimiento.Range("D16").Value = 28410.65
$wsCumplimiento.Range("E16").Value = 15855.59
$wsCumplimiento.Range("F16").Value = 0.6418130385594079

# Sheet "CUMPLIMIENTO MENSUAL": TOTAL row (VENTA / POR CUMPLIR / CUMPLIMIENTO)
$wsCumplimiento.Range("D19").Value = 37623.7
$wsCumplimiento.Range("E19").Value = 27754.29762291769
$wsCumplimiento.Range("F19").Value = 0.5754795400281782
